# Auto-generated data-driven update for Recommandations (sheet1) and Top_YTD (sheet2)
$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsTop  = $wb.Worksheets.Item("Top_YTD")

# --- Recommandations sheet: rows 2-48, columns A-G ---
$recoData = @(
    @(2, "SUCRIVOIRE", 0, 4, 3910, 960, "🟡 Observer", "➖ Neutre")
    ,@(3, "BRVM - SERVICES PUBLICS", 0, 8, 3429.52, 112.76, "🟡 Observer", "➖ Neutre")
    ,@(4, "SAFCA CI", 0, 4, 2775, 690, "🟡 Observer", "➖ Neutre")
    ,@(5, "CFAO MOTORS CI", 0, 4, 2700, 680, "🟡 Observer", "➖ Neutre")
    ,@(6, "BRVM - AUTRES SECTEURS", 0, 4, 2656.22, 661.58, "🟡 Observer", "➖ Neutre")
    ,@(7, "NEI-CEDA CI", 0, 4, 2385, 595, "🟡 Observer", "➖ Neutre")
    ,@(8, "UNIWAX CI", 0, 4, 2365, 590, "🟡 Observer", "➖ Neutre")
    ,@(9, "SETAO CI", 0, 4, 2245, 565, "🟡 Observer", "➖ Neutre")
    ,@(10, "AIR LIQUIDE CI", 0, 4, 2155, 535, "🟡 Observer", "➖ Neutre")
    ,@(11, "BRVM - DISTRIBUTION", 0, 4, 1478.55, 370.99, "🟡 Observer", "➖ Neutre")
    ,@(12, "BRVM - TRANSPORT", 0, 4, 1400.1, 351.25, "🟡 Observer", "➖ Neutre")
    ,@(13, "BRVM - AGRICULTURE", 0, 4, 1326.86, 331.52, "🟡 Observer", "➖ Neutre")
    ,@(14, "BRVM - INDUSTRIE", 0, 4, 772.34, 193.08, "🟡 Observer", "➖ Neutre")
    ,@(15, "BRVM-PRINCIPAL", 0, 4, 706.7, 175.95, "🟡 Observer", "➖ Neutre")
    ,@(16, "BRVM - CONSOMMATION DE BASE", 0, 4, 686.5599999999999, 172.4, "🟡 Observer", "➖ Neutre")
    ,@(17, "BRVM-PRESTIGE", 0, 4, 527.21, 131.55, "🟡 Observer", "➖ Neutre")
    ,@(18, "BRVM - INDUSTRIELS", 0, 4, 516.09, 125.96, "🟡 Observer", "➖ Neutre")
    ,@(19, "BRVM - FINANCES", 0, 4, 493.18, 123.03, "🟡 Observer", "➖ Neutre")
    ,@(20, "BRVM - SERVICES FINANCIERS", 0, 4, 484.69, 120.91, "🟡 Observer", "➖ Neutre")
    ,@(21, "BRVM - ENERGIE", 0, 4, 439.51, 110.38, "🟡 Observer", "➖ Neutre")
    ,@(22, "BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 425.46, 105.25, "🟡 Observer", "➖ Neutre")
    ,@(23, "BRVM - TELECOMMUNICATIONS", 0, 4, 387.33, 95.5, "🟡 Observer", "➖ Neutre")
    ,@(24, "UNILEVER CI (UNLC)", 3, 0, 22.45, 7.49, "🟢 Achat", "✅ Renforcer")
    ,@(25, "VIVO ENERGY CI (SHEC)", 2, 0, 5.37, 0.9399999999999999, "🟡 Observer", "➖ Neutre")
    ,@(26, "SETAO CI (STAC)", 1, 1, 4.82, -2.59, "🟡 Observer", "👀 À surveiller")
    ,@(27, "SODE CI (SDCC)", 2, 0, 4.1, 1.67, "🟡 Observer", "➖ Neutre")
    ,@(28, "TOTALENERGIES MARKETING CI (TTLC)", 1, 0, 3.39, 3.39, "🟡 Observer", "➖ Neutre")
    ,@(29, "BANK OF AFRICA SENEGAL (BOAS)", 1, 0, 3.38, 3.38, "🟡 Observer", "➖ Neutre")
    ,@(30, "BANK OF AFRICA ML (BOAM)", 1, 1, 3.29, 6.22, "🟡 Observer", "👀 À surveiller")
    ,@(31, "ORANGE COTE D'IVOIRE (ORAC)", 1, 0, 2.6, 2.6, "🟡 Observer", "➖ Neutre")
    ,@(32, "PALM CI (PALC)", 1, 1, 2.01, -2.74, "🟡 Observer", "👀 À surveiller")
    ,@(33, "SAFCA CI (SAFC)", 1, 0, 1.47, 1.47, "🟡 Observer", "➖ Neutre")
    ,@(34, "ONATEL BF (ONTBF)", 1, 0, 1.35, 1.35, "🟡 Observer", "➖ Neutre")
    ,@(35, "ECOBANK TRANS. INCORP. TG (ETIT)", 1, 1, 0.32, 5.88, "🟡 Observer", "👀 À surveiller")
    ,@(36, "TOTAL", 0, 4, 0, 0, "🟡 Observer", "➖ Neutre")
    ,@(37, "SONATEL SN (SNTS)", 1, 1, -0.62, 3.61, "🟡 Observer", "👀 À surveiller")
    ,@(38, "SERVAIR ABIDJAN CI (ABJC)", 1, 1, -0.93, 4.24, "🟡 Observer", "👀 À surveiller")
    ,@(39, "BANK OF AFRICA BF (BOABF)", 0, 1, -2.23, -2.23, "🟡 Observer", "➖ Neutre")
    ,@(40, "CFAO MOTORS CI (CFAC)", 0, 1, -2.94, -2.94, "🟡 Observer", "➖ Neutre")
    ,@(41, "SMB CI (SMBC)", 1, 2, -3.66, -3.24, "🟡 Observer", "👀 À surveiller")
    ,@(42, "TRACTAFRIC MOTORS CI (PRSC)", 0, 1, -3.85, -3.85, "🟡 Observer", "➖ Neutre")
    ,@(43, "BANK OF AFRICA BN (BOAB)", 0, 1, -3.95, -3.95, "🟡 Observer", "➖ Neutre")
    ,@(44, "BERNABE CI (BNBC)", 0, 1, -4.17, -4.17, "🟡 Observer", "➖ Neutre")
    ,@(45, "SOLIBRA CI (SLBC)", 0, 1, -6.67, -6.67, "🟡 Observer", "➖ Neutre")
    ,@(46, "AFRICA GLOBAL LOGISTICS CI (SDSC)", 0, 1, -7.42, -7.42, "🟡 Observer", "➖ Neutre")
    ,@(47, "BANK OF AFRICA NG (BOAN)", 1, 2, -7.89, -6.65, "🟡 Observer", "👀 À surveiller")
    ,@(48, "FILTISAC CI (FTSC)", 0, 3, -22.26, -7.4, "🔴 Vente", "⚠️ Risque de décrochage")
)

foreach ($row in $recoData) {
    $r = $row[0]
    $wsReco.Cells.Item($r, 1).Value = $row[1]
    $wsReco.Cells.Item($r, 2).Value = $row[2]
    $wsReco.Cells.Item($r, 3).Value = $row[3]
    $wsReco.Cells.Item($r, 4).Value = $row[4]
    $wsReco.Cells.Item($r, 5).Value = $row[5]
    $wsReco.Cells.Item($r, 6).Value = $row[6]
    $wsReco.Cells.Item($r, 7).Value = $row[7]
}

# --- Top_YTD sheet: rows 2-11, columns A-B ---
$topData = @(
    @(2, "BRVM - SERVICES PUBLICS", 10529595.09)
    ,@(3, "SUCRIVOIRE", 1347529.68)
    ,@(4, "SAFCA CI", 396827.6)
    ,@(5, "CFAO MOTORS CI", 360635.38)
    ,@(6, "BRVM - AUTRES SECTEURS", 340695.65)
    ,@(7, "NEI-CEDA CI", 234891.66)
    ,@(8, "UNIWAX CI", 228189.95)
    ,@(9, "SETAO CI", 190986.09)
    ,@(10, "AIR LIQUIDE CI", 166351.28)
    ,@(11, "BRVM - DISTRIBUTION", 48539.89)
)

foreach ($row in $topData) {
    $r = $row[0]
    $wsTop.Cells.Item($r, 1).Value = $row[1]
    $wsTop.Cells.Item($r, 2).Value = $row[2]
}

Write-Host "Update complete: Recommandations rows=$($recoData.Length), Top_YTD rows=$($topData.Length)"
